$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @("H19", 415.66666),
    @("I19", 618.4),
    @("K19", 618.4),
    @("M19", -443.4),
    @("H53", 423.55554),
    @("I53", 398),
    @("J53", 474.66666),
    @("K53", 398),
    @("L53", 474.66666),
    @("M53", 239),
    @("N53", -1748.66666),
    @("H86", 5152.4443),
    @("I86", 5895.6665),
    @("J86", 3666),
    @("K86", 5895.6665),
    @("L86", 3666),
    @("M86", -4772.6665),
    @("N86", -5912),
    @("H89", 5152.4443),
    @("I89", 5895.6665),
    @("J89", 3666),
    @("K89", 29478.3325),
    @("L89", 18330),
    @("M89", -23862.3325),
    @("N89", -29562),
    @("H112", 2209.261),
    @("J112", 2373.3572),
    @("L112", 7120.071599999999),
    @("N112", -9336.071599999999),
    @("H116", 13492.333),
    @("I116", 5080.625),
    @("J116", 20221.7),
    @("K116", 5080.625),
    @("L116", 20221.7),
    @("M116", -1638.625),
    @("N116", -27105.7),
    @("H141", 3931.889),
    @("I141", 3327.8),
    @("J141", 4687),
    @("K141", 9983.400000000001),
    @("L141", 14061),
    @("M141", -4803.400000000001),
    @("N141", -24421)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @("H2", 719.8),
    @("I2", 719.8),
    @("K2", 719.8),
    @("M2", -606.8),
    @("H32", 5313.1387),
    @("I32", 4548.373),
    @("J32", 12833.333),
    @("K32", 4548.373),
    @("L32", 12833.333),
    @("M32", -4261.373),
    @("N32", -13407.333),
    @("H38", 5777.5),
    @("I38", 1534.5),
    @("J38", 10020.5),
    @("K38", 1534.5),
    @("L38", 10020.5),
    @("M38", -1067.5),
    @("N38", -10954.5),
    @("H41", 2653.0908),
    @("I41", 2168.4),
    @("K41", 2168.4),
    @("M41", -1754.4),
    @("H63", 1872.6364),
    @("I63", 1872.6364),
    @("K63", 1872.6364),
    @("M63", -1186.6364),
    @("H66", 1872.6364),
    @("I66", 1872.6364),
    @("K66", 9363.182000000001),
    @("M66", -5931.182000000001),
    @("H74", 4654.983),
    @("I74", 964.3555),
    @("J74", 17430.23),
    @("K74", 964.3555),
    @("L74", 17430.23),
    @("M74", -90.35550000000001),
    @("N74", -19178.23),
    @("H77", 4654.983),
    @("I77", 964.3555),
    @("J77", 17430.23),
    @("K77", 4821.7775),
    @("L77", 87151.14999999999),
    @("M77", -453.7775000000001),
    @("N77", -95887.14999999999),
    @("H97", 40001504),
    @("I97", 1200.826),
    @("K97", 1200.826),
    @("M97", -704.826),
    @("H102", 3359.125),
    @("I102", 3359.125),
    @("K102", 3359.125),
    @("M102", -1737.125),
    @("H110", 1004.4167),
    @("I110", 913.9091),
    @("K110", 913.9091),
    @("M110", 1131.0909),
    @("H116", 719.8),
    @("I116", 719.8),
    @("K116", 719.8),
    @("M116", 1574.2)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @("H3", 719.8),
    @("I3", 719.8),
    @("K3", 719.8),
    @("M3", -605.8),
    @("H86", 5000),
    @("J86", 5000),
    @("L86", 5000),
    @("N86", -7246),
    @("H89", 5000),
    @("J89", 5000),
    @("L89", 25000),
    @("N89", -36232),
    @("H107", 824),
    @("I107", 794.6667),
    @("J107", 1000),
    @("K107", 794.6667),
    @("L107", 1000),
    @("M107", 1125.3333),
    @("N107", -4840)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @("H16", 1658.3334),
    @("I16", 1658.3334),
    @("K16", 1658.3334),
    @("M16", -1371.3334),
    @("H25", 1179.3572),
    @("I25", 1387.3636),
    @("K25", 1387.3636),
    @("M25", -1213.3636),
    @("H31", 3079.0667),
    @("I31", 3749.4285),
    @("K31", 3749.4285),
    @("M31", -3454.4285),
    @("H34", 3079.0667),
    @("I34", 3749.4285),
    @("K34", 3749.4285),
    @("M34", -3547.4285),
    @("H86", 13327.091),
    @("J86", 4816.3335),
    @("L86", 4816.3335),
    @("N86", -7062.3335),
    @("H89", 13327.091),
    @("J89", 4816.3335),
    @("L89", 24081.6675),
    @("N89", -35313.6675),
    @("H113", 1658.3334),
    @("I113", 1658.3334),
    @("K113", 1658.3334),
    @("M113", 511.6666),
    @("H134", 3309.0908),
    @("I134", 2764.7058),
    @("J134", 5160),
    @("K134", 8294.117400000001),
    @("L134", 15480),
    @("M134", -5759.117400000001),
    @("N134", -20550)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    @("H131", 4631.5454),
    @("J131", 5011.55),
    @("L131", 15034.65),
    @("N131", -25114.65)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    @("H70", 23876),
    @("I70", 12193.059),
    @("K70", 12193.059),
    @("M70", -11923.059),
    @("H73", 23876),
    @("I73", 12193.059),
    @("K73", 12193.059),
    @("M73", -11257.059),
    @("H80", 0),
    @("I80", 0),
    @("K80", 0),
    @("M80", $null),
    @("H83", 0),
    @("I83", 0),
    @("K83", 0),
    @("M83", $null),
    @("H122", 4081.4375),
    @("I122", 3730.4),
    @("J122", 4666.5),
    @("K122", 11191.2),
    @("L122", 13999.5),
    @("M122", -8741.200000000001),
    @("N122", -18899.5),
    @("H123", 54999),
    @("J123", 54999),
    @("L123", 54999),
    @("N123", -59899),
    @("H132", 9183.111000000001),
    @("I132", 10880.464),
    @("K132", 32641.392),
    @("M132", -30111.392)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @("H29", 0),
    @("J29", 0),
    @("L29", 0),
    @("N29", $null),
    @("H61", 9195.799999999999),
    @("I61", 3333.7),
    @("J61", 20920),
    @("K61", 3333.7),
    @("L61", 20920),
    @("M61", -3131.7),
    @("N61", -21324),
    @("H93", 4240.222),
    @("I93", 1488.3334),
    @("K93", 1488.3334),
    @("M93", -240.3334),
    @("H113", 9195.799999999999),
    @("I113", 3333.7),
    @("J113", 20920),
    @("K113", 3333.7),
    @("L113", 20920),
    @("M113", -1163.7),
    @("N113", -25260),
    @("H141", 50000),
    @("J141", 50000),
    @("L141", 50000),
    @("N141", -60360)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @("H96", 2143.5),
    @("I96", 1735.4286),
    @("J96", 5000),
    @("K96", 1735.4286),
    @("L96", 5000),
    @("M96", -362.4286),
    @("N96", -7746),
    @("H122", 87673.62),
    @("I122", 767.625),
    @("J122", 226723.2),
    @("K122", 2302.875),
    @("L122", 680169.6000000001),
    @("M122", 147.125),
    @("N122", -685069.6000000001)
)
foreach ($edit in $edits) {
    $addr = $edit[0]
    $val = $edit[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}
